$wb = $excel.ActiveWorkbook

# Update the id values in Sheet1 (D7:D12) - increment each by its row offset
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("D7").Value = 1033
$ws1.Range("D8").Value = 1034
$ws1.Range("D9").Value = 1035
$ws1.Range("D10").Value = 1036
$ws1.Range("D11").Value = 1037
$ws1.Range("D12").Value = 1038

# Make Sheet1 the active sheet/tab (was "thirdSheet" before)
$ws1.Activate()
